$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 105, shifting rows 105-225 down to 106-226
$ws.Rows.Item(105).Insert()

# Populate the new row 105 with the data from the diff
$ws.Cells.Item(105, 1).Value = 5
$ws.Cells.Item(105, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(105, 3).Value = "Maule"
$ws.Cells.Item(105, 4).Value = 44546
$ws.Cells.Item(105, 5).Value = 7
$ws.Cells.Item(105, 6).Value = 100114014
$ws.Cells.Item(105, 7).Value = "Betarraga"
$ws.Cells.Item(105, 8).Value = "Sin especificar"
$ws.Cells.Item(105, 9).Value = "Primera"
$ws.Cells.Item(105, 10).Value = 5000
$ws.Cells.Item(105, 11).Value = 500
$ws.Cells.Item(105, 12).Value = 500
$ws.Cells.Item(105, 13).Value = 500
$ws.Cells.Item(105, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(105, 15).Value = "Región del Maule"
$ws.Cells.Item(105, 16).Value = 100
$ws.Cells.Item(105, 17).Value = 5
$ws.Cells.Item(105, 18).Value = "Hortaliza"
